$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A1").Value = 1.24846339225769
$ws.Range("B1").Value = 1.136787891387939
$ws.Range("C1").Value = 5.170534610748291
$ws.Range("D1").Value = 1.550818562507629
$ws.Range("E1").Value = 1.017218828201294
